$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update team colors
$ws.Range("C15").Value = "#F1AA11"   # Los Angeles Lakers
$ws.Range("C29").Value = "#753BBD"   # Toronto Raptors
$ws.Range("C12").Value = "#BA0C2F"   # Houston Rockets

# Update the active cell selection to match the saved view
$ws.Range("C13").Select()
